# Applies per-row price/volume/coin updates to Sheet1 as described in the commit diff.
# D-column values that parse as plain numbers need NumberFormat "@" (Text) set
# BEFORE the value is assigned, otherwise Excel COM auto-converts the numeric-looking
# string (e.g. "486.30") into a float and the trailing zero / exact text is lost.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "68.267.81"
$ws.Cells.Item(2, 5).Value = "  +0.12%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "3.917.04"
$ws.Cells.Item(3, 5).Value = "  -1.37%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  +0.19%  "

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "486.30"
$ws.Cells.Item(5, 5).Value = "  +0.03%  "

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "145.93"
$ws.Cells.Item(6, 5).Value = "  -2.37%  "

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.623"
$ws.Cells.Item(7, 5).Value = "  -0.96%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  -0.07%  "

# Row 9
$ws.Cells.Item(9, 5).Value = "  -0.49%  "

# Row 10
$ws.Cells.Item(10, 5).Value = "  -2.58%  "

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "0.0000349"
$ws.Cells.Item(11, 5).Value = "  -6.06%  "

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "43.10"
$ws.Cells.Item(12, 5).Value = "  -1.57%  "

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "10.71"
$ws.Cells.Item(13, 5).Value = "  +1.85%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "4.545.59"
$ws.Cells.Item(14, 5).Value = "  -1.11%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "3.911.15"
$ws.Cells.Item(15, 5).Value = "  -1.28%  "

# Row 16
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "14.32"
$ws.Cells.Item(16, 5).Value = "  -3.99%  "

# Row 17
$ws.Cells.Item(17, 5).Value = "  -0.78%  "

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "20.09"
$ws.Cells.Item(18, 5).Value = "  +0.57%  "

# Row 19
$ws.Cells.Item(19, 5).Value = "  -0.90%  "

# Row 20
$ws.Cells.Item(20, 4).Value = "68.352.61"

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "431.80"
$ws.Cells.Item(21, 5).Value = "  -1.23%  "

# Row 22
$ws.Cells.Item(22, 5).Value = "  +3.42%  "

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "15.12"
$ws.Cells.Item(23, 5).Value = "  +4.98%  "

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "88.49"
$ws.Cells.Item(24, 5).Value = "  +0.27%  "

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "11.73"
$ws.Cells.Item(25, 5).Value = "  +21.27%  "

# Row 26
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "3.71"
$ws.Cells.Item(26, 5).Value = "  +1.40%  "

# Row 27
$ws.Cells.Item(27, 5).Value = "  +10.77%  "

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "37.82"
$ws.Cells.Item(28, 5).Value = "  -3.21%  "

# Row 29
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "5.67"
$ws.Cells.Item(29, 5).Value = "  -1.32%  "

# Row 30
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "719.02"
$ws.Cells.Item(30, 5).Value = "  -1.74%  "

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "13.77"
$ws.Cells.Item(31, 5).Value = "  +3.33%  "

# Row 32
$ws.Cells.Item(32, 5).Value = "  +0.96%  "

# Row 33
$ws.Cells.Item(33, 5).Value = "  +3.16%  "

# Row 34
$ws.Cells.Item(34, 4).Value = "0.0₃0913"
$ws.Cells.Item(34, 5).Value = "  +1.77%  "

# Row 35
$ws.Cells.Item(35, 5).Value = "  +15.24%  "

# Row 36
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "41.79"
$ws.Cells.Item(36, 5).Value = "  -0.45%  "

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "60.86"
$ws.Cells.Item(37, 5).Value = "  +0.72%  "

# Row 38
$ws.Cells.Item(38, 2).Value = "TheGraph"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.397"
$ws.Cells.Item(38, 5).Value = "  +18.52%  "

# Row 39
$ws.Cells.Item(39, 2).Value = "Fetch.AI"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "3.01"
$ws.Cells.Item(39, 5).Value = "  +19.85%  "

# Row 40
$ws.Cells.Item(40, 2).Value = "Dai"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.998"
$ws.Cells.Item(40, 5).Value = "  +0.03%  "

# Row 41
$ws.Cells.Item(41, 2).Value = "Kaspa"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.145"
$ws.Cells.Item(41, 5).Value = "  -4.41%  "

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.0492"
$ws.Cells.Item(42, 5).Value = "  +3.60%  "

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "3.13"
$ws.Cells.Item(43, 5).Value = "  +2.29%  "

# Row 44
$ws.Cells.Item(44, 5).Value = "  +5.08%  "

# Row 45
$ws.Cells.Item(45, 5).Value = "  -0.19%  "

# Row 46
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "3.36"
$ws.Cells.Item(46, 5).Value = "  +4.01%  "

# Row 47
$ws.Cells.Item(47, 5).Value = "  +0.00%  "

# Row 48
$ws.Cells.Item(48, 5).Value = "  -0.41%  "

# Row 49
$ws.Cells.Item(49, 5).Value = "  -5.35%  "

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "144.93"
$ws.Cells.Item(50, 5).Value = "  -2.24%  "

# Row 51
$ws.Cells.Item(51, 4).Value = "0.0₆0341"
$ws.Cells.Item(51, 5).Value = "  +27.77%  "

